$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.999.24"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "2.534.33"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "317.57"
$ws.Range("E5").Value = "  +4.39%  "
$ws.Range("D6").Value = "95.94"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").Value = "36.31"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "7.63"
$ws.Range("E12").Value = "  -0.59%  "
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "2.924.33"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.526.12"
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "15.40"
$ws.Range("E16").Value = "  +2.39%  "
$ws.Range("D17").Value = "0.852"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "43.048.13"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").Value = "6.68"
$ws.Range("E20").Value = "  +3.63%  "
$ws.Range("D21").Value = "0.0₃0969"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").Value = "70.30"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("D23").Value = "252.31"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("D24").Value = "2.97"
$ws.Range("E24").Value = "  +1.76%  "
$ws.Range("D25").Value = "2.02"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "2.43"
$ws.Range("E28").Value = "  +5.12%  "
$ws.Range("D29").Value = "40.06"
$ws.Range("E29").Value = "  +5.50%  "
$ws.Range("D30").Value = "10.31"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").Value = "6.05"
$ws.Range("E31").Value = "  +1.66%  "
$ws.Range("D32").Value = "154.71"
$ws.Range("E32").Value = "  -1.45%  "
$ws.Range("E33").Value = "  +2.96%  "
$ws.Range("D34").Value = "3.31"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").Value = "0.0793"
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("D36").Value = "18.82"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "23.77"
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("E41").Value = "  +13.84%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "3.83"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "0.0306"
$ws.Range("E43").Value = "  +2.09%  "
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").Value = "3.30"
$ws.Range("E45").Value = "  -2.61%  "
$ws.Range("D46").Value = "2.023.75"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("D47").Value = "86.00"
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("D48").Value = "8.83"
$ws.Range("E48").Value = "  -1.70%  "
$ws.Range("D49").Value = "2.779.55"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").Value = "74.01"
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").Value = "102.88"
$ws.Range("E51").Value = "  +1.36%  "
